# Sample Project / Main.xlsx - "Rules" sheet
# Cell B11 (the rule name for what used to be rule "R40") becomes the
# literal text "1" while keeping its existing cell style/format.
#
# A plain  $ws.Range("B11").Value = "1"  would be auto-coerced to the
# *number* 1 by Excel's type inference (since "1" parses as numeric),
# and forcing a Text number format first marks the cell with a
# quote-prefix / new style id instead of reusing the original one.
# Entering it as a formula that evaluates to the text "1" and then
# collapsing that formula down to its literal value via Paste Special
# (values only) reproduces exactly what a normal "type text into the
# cell" edit looks like on disk: a shared-string cell (t="s") that
# keeps its original style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cell = $ws.Range("B11")
$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false
